$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for Jengibre was added to the source feed.
# It belongs right after the current first data block for this product,
# so insert a fresh row at 322 (this shifts the existing rows 322:349
# down to 323:350) and populate it with the new record's values.
$ws.Rows.Item(322).Insert()

$ws.Cells.Item(322, 1).Value = 10
$ws.Cells.Item(322, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(322, 3).Value = "La Araucanía"
$ws.Cells.Item(322, 4).Value = 45166
$ws.Cells.Item(322, 5).Value = 9
$ws.Cells.Item(322, 6).Value = 100114007
$ws.Cells.Item(322, 7).Value = "Jengibre"
$ws.Cells.Item(322, 8).Value = "Sin especificar"
$ws.Cells.Item(322, 9).Value = "Primera"
$ws.Cells.Item(322, 10).Value = 180
$ws.Cells.Item(322, 11).Value = 22000
$ws.Cells.Item(322, 12).Value = 24000
$ws.Cells.Item(322, 13).Value = 23111
$ws.Cells.Item(322, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(322, 15).Value = "Perú"
$ws.Cells.Item(322, 16).Value = 1778
$ws.Cells.Item(322, 17).Value = 13
$ws.Cells.Item(322, 18).Value = "Hortaliza"
